$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 18
$ws.Range("A18").Value = 7
$ws.Range("B18").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C18").Value = "Ñuble"
$ws.Range("D18").Value = 44664
$ws.Range("E18").Value = 16
$ws.Range("F18").Value = 100112040
$ws.Range("G18").Value = "Cilantro"
$ws.Range("H18").Value = "Sin especificar"
$ws.Range("I18").Value = "Primera"
$ws.Range("J18").Value = 200
$ws.Range("K18").Value = 550
$ws.Range("L18").Value = 600
$ws.Range("M18").Value = 575
$ws.Range("N18").Value = "`$/atado 0,5 a 1 kilo"
$ws.Range("O18").Value = "Provincia de Diguillín"
$ws.Range("P18").Value = 575
$ws.Range("Q18").Value = 1
$ws.Range("R18").Value = "Hortaliza"

# Row 19
$ws.Range("A19").Value = 7
$ws.Range("B19").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C19").Value = "Ñuble"
$ws.Range("D19").Value = 44623
$ws.Range("E19").Value = 16
$ws.Range("F19").Value = 100112040
$ws.Range("G19").Value = "Cilantro"
$ws.Range("H19").Value = "Sin especificar"
$ws.Range("I19").Value = "Primera"
$ws.Range("J19").Value = 120
$ws.Range("K19").Value = 550
$ws.Range("L19").Value = 600
$ws.Range("M19").Value = 575
$ws.Range("N19").Value = "`$/atado 0,5 a 1 kilo"
$ws.Range("O19").Value = "Provincia de Diguillín"
$ws.Range("P19").Value = 575
$ws.Range("Q19").Value = 1
$ws.Range("R19").Value = "Hortaliza"

# Row 20
$ws.Range("A20").Value = 7
$ws.Range("B20").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C20").Value = "Ñuble"
$ws.Range("D20").Value = 44656
$ws.Range("E20").Value = 16
$ws.Range("F20").Value = 100112040
$ws.Range("G20").Value = "Cilantro"
$ws.Range("H20").Value = "Sin especificar"
$ws.Range("I20").Value = "Primera"
$ws.Range("J20").Value = 200
$ws.Range("K20").Value = 600
$ws.Range("L20").Value = 650
$ws.Range("M20").Value = 625
$ws.Range("N20").Value = "`$/atado 0,5 a 1 kilo"
$ws.Range("O20").Value = "Provincia de Diguillín"
$ws.Range("P20").Value = 625
$ws.Range("Q20").Value = 1
$ws.Range("R20").Value = "Hortaliza"

# Row 21
$ws.Range("A21").Value = 7
$ws.Range("B21").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C21").Value = "Ñuble"
$ws.Range("D21").Value = 44630
$ws.Range("E21").Value = 16
$ws.Range("F21").Value = 100112040
$ws.Range("G21").Value = "Cilantro"
$ws.Range("H21").Value = "Sin especificar"
$ws.Range("I21").Value = "Primera"
$ws.Range("J21").Value = 200
$ws.Range("K21").Value = 550
$ws.Range("L21").Value = 600
$ws.Range("M21").Value = 575
$ws.Range("N21").Value = "`$/atado 0,5 a 1 kilo"
$ws.Range("O21").Value = "Provincia de Diguillín"
$ws.Range("P21").Value = 575
$ws.Range("Q21").Value = 1
$ws.Range("R21").Value = "Hortaliza"

# Row 22
$ws.Range("A22").Value = 7
$ws.Range("B22").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C22").Value = "Ñuble"
$ws.Range("D22").Value = 44267
$ws.Range("E22").Value = 16
$ws.Range("F22").Value = 100112040
$ws.Range("G22").Value = "Cilantro"
$ws.Range("H22").Value = "Sin especificar"
$ws.Range("I22").Value = "Primera"
$ws.Range("J22").Value = 150
$ws.Range("K22").Value = 1800
$ws.Range("L22").Value = 2000
$ws.Range("M22").Value = 1913
$ws.Range("N22").Value = "`$/atado 0,5 a 1 kilo"
$ws.Range("O22").Value = "Provincia de Diguillín"
$ws.Range("P22").Value = 1913
$ws.Range("Q22").Value = 1
$ws.Range("R22").Value = "Hortaliza"

# Row 23
$ws.Range("A23").Value = 7
$ws.Range("B23").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C23").Value = "Ñuble"
$ws.Range("D23").Value = 44649
$ws.Range("E23").Value = 16
$ws.Range("F23").Value = 100112040
$ws.Range("G23").Value = "Cilantro"
$ws.Range("H23").Value = "Sin especificar"
$ws.Range("I23").Value = "Primera"
$ws.Range("J23").Value = 200
$ws.Range("K23").Value = 600
$ws.Range("L23").Value = 650
$ws.Range("M23").Value = 625
$ws.Range("N23").Value = "`$/atado 0,5 a 1 kilo"
$ws.Range("O23").Value = "Provincia de Diguillín"
$ws.Range("P23").Value = 625
$ws.Range("Q23").Value = 1
$ws.Range("R23").Value = "Hortaliza"

# Row 24
$ws.Range("A24").Value = 7
$ws.Range("B24").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C24").Value = "Ñuble"
$ws.Range("D24").Value = 44272
$ws.Range("E24").Value = 16
$ws.Range("F24").Value = 100112040
$ws.Range("G24").Value = "Cilantro"
$ws.Range("H24").Value = "Sin especificar"
$ws.Range("I24").Value = "Primera"
$ws.Range("J24").Value = 150
$ws.Range("K24").Value = 1800
$ws.Range("L24").Value = 2000
$ws.Range("M24").Value = 1893
$ws.Range("N24").Value = "`$/atado 0,5 a 1 kilo"
$ws.Range("O24").Value = "Provincia de Diguillín"
$ws.Range("P24").Value = 1893
$ws.Range("Q24").Value = 1
$ws.Range("R24").Value = "Hortaliza"

# Row 25
$ws.Range("A25").Value = 7
$ws.Range("B25").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C25").Value = "Ñuble"
$ws.Range("D25").Value = 44533
$ws.Range("E25").Value = 16
$ws.Range("F25").Value = 100112040
$ws.Range("G25").Value = "Cilantro"
$ws.Range("H25").Value = "Sin especificar"
$ws.Range("I25").Value = "Primera"
$ws.Range("J25").Value = 100
$ws.Range("K25").Value = 2000
$ws.Range("L25").Value = 2200
$ws.Range("M25").Value = 2100
$ws.Range("N25").Value = "`$/atado 0,5 a 1 kilo"
$ws.Range("O25").Value = "Provincia de Diguillín"
$ws.Range("P25").Value = 2100
$ws.Range("Q25").Value = 1
$ws.Range("R25").Value = "Hortaliza"

# Row 26
$ws.Range("A26").Value = 7
$ws.Range("B26").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C26").Value = "Ñuble"
$ws.Range("D26").Value = 44608
$ws.Range("E26").Value = 16
$ws.Range("F26").Value = 100112040
$ws.Range("G26").Value = "Cilantro"
$ws.Range("H26").Value = "Sin especificar"
$ws.Range("I26").Value = "Primera"
$ws.Range("J26").Value = 160
$ws.Range("K26").Value = 550
$ws.Range("L26").Value = 600
$ws.Range("M26").Value = 575
$ws.Range("N26").Value = "`$/atado 0,5 a 1 kilo"
$ws.Range("O26").Value = "Provincia de Diguillín"
$ws.Range("P26").Value = 575
$ws.Range("Q26").Value = 1
$ws.Range("R26").Value = "Hortaliza"

# Row 27
$ws.Range("A27").Value = 7
$ws.Range("B27").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C27").Value = "Ñuble"
$ws.Range("D27").Value = 44644
$ws.Range("E27").Value = 16
$ws.Range("F27").Value = 100112040
$ws.Range("G27").Value = "Cilantro"
$ws.Range("H27").Value = "Sin especificar"
$ws.Range("I27").Value = "Primera"
$ws.Range("J27").Value = 100
$ws.Range("K27").Value = 550
$ws.Range("L27").Value = 600
$ws.Range("M27").Value = 575
$ws.Range("N27").Value = "`$/atado 0,5 a 1 kilo"
$ws.Range("O27").Value = "Provincia de Diguillín"
$ws.Range("P27").Value = 575
$ws.Range("Q27").Value = 1
$ws.Range("R27").Value = "Hortaliza"

# Row 28
$ws.Range("A28").Value = 7
$ws.Range("B28").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C28").Value = "Ñuble"
$ws.Range("D28").Value = 44160
$ws.Range("E28").Value = 16
$ws.Range("F28").Value = 100112040
$ws.Range("G28").Value = "Cilantro"
$ws.Range("H28").Value = "Sin especificar"
$ws.Range("I28").Value = "Primera"
$ws.Range("J28").Value = 190
$ws.Range("K28").Value = 1300
$ws.Range("L28").Value = 1500
$ws.Range("M28").Value = 1395
$ws.Range("N28").Value = "`$/atado 1 a 1,5 kilos"
$ws.Range("O28").Value = "Provincia de Diguillín"
$ws.Range("P28").Value = 930
$ws.Range("Q28").Value = 1.5
$ws.Range("R28").Value = "Hortaliza"

# Row 29
$ws.Range("A29").Value = 7
$ws.Range("B29").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C29").Value = "Ñuble"
$ws.Range("D29").Value = 44159
$ws.Range("E29").Value = 16
$ws.Range("F29").Value = 100112040
$ws.Range("G29").Value = "Cilantro"
$ws.Range("H29").Value = "Sin especificar"
$ws.Range("I29").Value = "Primera"
$ws.Range("J29").Value = 55
$ws.Range("K29").Value = 7000
$ws.Range("L29").Value = 8000
$ws.Range("M29").Value = 7455
$ws.Range("N29").Value = "`$/caja 36 atados"
$ws.Range("O29").Value = "Región Metropolitana"
$ws.Range("P29").Value = 207
$ws.Range("Q29").Value = 36
$ws.Range("R29").Value = "Hortaliza"

# Row 30
$ws.Range("A30").Value = 7
$ws.Range("B30").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C30").Value = "Ñuble"
$ws.Range("D30").Value = 44651
$ws.Range("E30").Value = 16
$ws.Range("F30").Value = 100112040
$ws.Range("G30").Value = "Cilantro"
$ws.Range("H30").Value = "Sin especificar"
$ws.Range("I30").Value = "Primera"
$ws.Range("J30").Value = 200
$ws.Range("K30").Value = 600
$ws.Range("L30").Value = 650
$ws.Range("M30").Value = 625
$ws.Range("N30").Value = "`$/atado 0,5 a 1 kilo"
$ws.Range("O30").Value = "Provincia de Diguillín"
$ws.Range("P30").Value = 625
$ws.Range("Q30").Value = 1
$ws.Range("R30").Value = "Hortaliza"

# Row 31
$ws.Range("A31").Value = 7
$ws.Range("B31").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C31").Value = "Ñuble"
$ws.Range("D31").Value = 44628
$ws.Range("E31").Value = 16
$ws.Range("F31").Value = 100112040
$ws.Range("G31").Value = "Cilantro"
$ws.Range("H31").Value = "Sin especificar"
$ws.Range("I31").Value = "Primera"
$ws.Range("J31").Value = 240
$ws.Range("K31").Value = 550
$ws.Range("L31").Value = 600
$ws.Range("M31").Value = 575
$ws.Range("N31").Value = "`$/atado 0,5 a 1 kilo"
$ws.Range("O31").Value = "Provincia de Diguillín"
$ws.Range("P31").Value = 575
$ws.Range("Q31").Value = 1
$ws.Range("R31").Value = "Hortaliza"

# Row 32
$ws.Range("A32").Value = 7
$ws.Range("B32").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C32").Value = "Ñuble"
$ws.Range("D32").Value = 44265
$ws.Range("E32").Value = 16
$ws.Range("F32").Value = 100112040
$ws.Range("G32").Value = "Cilantro"
$ws.Range("H32").Value = "Sin especificar"
$ws.Range("I32").Value = "Primera"
$ws.Range("J32").Value = 220
$ws.Range("K32").Value = 1800
$ws.Range("L32").Value = 2000
$ws.Range("M32").Value = 1909
$ws.Range("N32").Value = "`$/atado 0,5 a 1 kilo"
$ws.Range("O32").Value = "Provincia de Diguillín"
$ws.Range("P32").Value = 1909
$ws.Range("Q32").Value = 1
$ws.Range("R32").Value = "Hortaliza"

# Row 33
$ws.Range("A33").Value = 7
$ws.Range("B33").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C33").Value = "Ñuble"
$ws.Range("D33").Value = 44635
$ws.Range("E33").Value = 16
$ws.Range("F33").Value = 100112040
$ws.Range("G33").Value = "Cilantro"
$ws.Range("H33").Value = "Sin especificar"
$ws.Range("I33").Value = "Primera"
$ws.Range("J33").Value = 120
$ws.Range("K33").Value = 550
$ws.Range("L33").Value = 600
$ws.Range("M33").Value = 575
$ws.Range("N33").Value = "`$/atado 0,5 a 1 kilo"
$ws.Range("O33").Value = "Provincia de Diguillín"
$ws.Range("P33").Value = 575
$ws.Range("Q33").Value = 1
$ws.Range("R33").Value = "Hortaliza"

# Row 34
$ws.Range("A34").Value = 7
$ws.Range("B34").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C34").Value = "Ñuble"
$ws.Range("D34").Value = 44271
$ws.Range("E34").Value = 16
$ws.Range("F34").Value = 100112040
$ws.Range("G34").Value = "Cilantro"
$ws.Range("H34").Value = "Sin especificar"
$ws.Range("I34").Value = "Primera"
$ws.Range("J34").Value = 200
$ws.Range("K34").Value = 1800
$ws.Range("L34").Value = 2000
$ws.Range("M34").Value = 1920
$ws.Range("N34").Value = "`$/atado 0,5 a 1 kilo"
$ws.Range("O34").Value = "Provincia de Diguillín"
$ws.Range("P34").Value = 1920
$ws.Range("Q34").Value = 1
$ws.Range("R34").Value = "Hortaliza"

# Row 35
$ws.Range("A35").Value = 7
$ws.Range("B35").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C35").Value = "Ñuble"
$ws.Range("D35").Value = 44663
$ws.Range("E35").Value = 16
$ws.Range("F35").Value = 100112040
$ws.Range("G35").Value = "Cilantro"
$ws.Range("H35").Value = "Sin especificar"
$ws.Range("I35").Value = "Primera"
$ws.Range("J35").Value = 200
$ws.Range("K35").Value = 550
$ws.Range("L35").Value = 600
$ws.Range("M35").Value = 575
$ws.Range("N35").Value = "`$/atado 0,5 a 1 kilo"
$ws.Range("O35").Value = "Provincia de Diguillín"
$ws.Range("P35").Value = 575
$ws.Range("Q35").Value = 1
$ws.Range("R35").Value = "Hortaliza"

# Row 36
$ws.Range("A36").Value = 7
$ws.Range("B36").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C36").Value = "Ñuble"
$ws.Range("D36").Value = 44524
$ws.Range("E36").Value = 16
$ws.Range("F36").Value = 100112040
$ws.Range("G36").Value = "Cilantro"
$ws.Range("H36").Value = "Sin especificar"
$ws.Range("I36").Value = "Primera"
$ws.Range("J36").Value = 80
$ws.Range("K36").Value = 2000
$ws.Range("L36").Value = 2000
$ws.Range("M36").Value = 2000
$ws.Range("N36").Value = "`$/atado 0,5 a 1 kilo"
$ws.Range("O36").Value = "Provincia de Diguillín"
$ws.Range("P36").Value = 2000
$ws.Range("Q36").Value = 1
$ws.Range("R36").Value = "Hortaliza"

# Row 37
$ws.Range("A37").Value = 7
$ws.Range("B37").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C37").Value = "Ñuble"
$ws.Range("D37").Value = 44620
$ws.Range("E37").Value = 16
$ws.Range("F37").Value = 100112040
$ws.Range("G37").Value = "Cilantro"
$ws.Range("H37").Value = "Sin especificar"
$ws.Range("I37").Value = "Primera"
$ws.Range("J37").Value = 160
$ws.Range("K37").Value = 550
$ws.Range("L37").Value = 600
$ws.Range("M37").Value = 575
$ws.Range("N37").Value = "`$/atado 0,5 a 1 kilo"
$ws.Range("O37").Value = "Provincia de Diguillín"
$ws.Range("P37").Value = 575
$ws.Range("Q37").Value = 1
$ws.Range("R37").Value = "Hortaliza"

# Row 38
$ws.Range("A38").Value = 7
$ws.Range("B38").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C38").Value = "Ñuble"
$ws.Range("D38").Value = 44211
$ws.Range("E38").Value = 16
$ws.Range("F38").Value = 100112040
$ws.Range("G38").Value = "Cilantro"
$ws.Range("H38").Value = "Sin especificar"
$ws.Range("I38").Value = "Primera"
$ws.Range("J38").Value = 120
$ws.Range("K38").Value = 1800
$ws.Range("L38").Value = 2000
$ws.Range("M38").Value = 1883
$ws.Range("N38").Value = "`$/atado 0,5 a 1 kilo"
$ws.Range("O38").Value = "Provincia de Diguillín"
$ws.Range("P38").Value = 1883
$ws.Range("Q38").Value = 1
$ws.Range("R38").Value = "Hortaliza"

# Row 39
$ws.Range("A39").Value = 7
$ws.Range("B39").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C39").Value = "Ñuble"
$ws.Range("D39").Value = 44266
$ws.Range("D39").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E39").Value = 16
$ws.Range("F39").Value = 100112040
$ws.Range("G39").Value = "Cilantro"
$ws.Range("H39").Value = "Sin especificar"
$ws.Range("I39").Value = "Primera"
$ws.Range("J39").Value = 150
$ws.Range("K39").Value = 1800
$ws.Range("L39").Value = 2000
$ws.Range("M39").Value = 1913
$ws.Range("N39").Value = "`$/atado 0,5 a 1 kilo"
$ws.Range("O39").Value = "Provincia de Diguillín"
$ws.Range("P39").Value = 1913
$ws.Range("Q39").Value = 1
$ws.Range("R39").Value = "Hortaliza"

# Row 40
$ws.Range("A40").Value = 7
$ws.Range("B40").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C40").Value = "Ñuble"
$ws.Range("D40").Value = 44609
$ws.Range("D40").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E40").Value = 16
$ws.Range("F40").Value = 100112040
$ws.Range("G40").Value = "Cilantro"
$ws.Range("H40").Value = "Sin especificar"
$ws.Range("I40").Value = "Primera"
$ws.Range("J40").Value = 60
$ws.Range("K40").Value = 550
$ws.Range("L40").Value = 600
$ws.Range("M40").Value = 575
$ws.Range("N40").Value = "`$/atado 0,5 a 1 kilo"
$ws.Range("O40").Value = "Provincia de Diguillín"
$ws.Range("P40").Value = 575
$ws.Range("Q40").Value = 1
$ws.Range("R40").Value = "Hortaliza"
